$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header area
# ---------------------------------------------------------------------------

# C2 rich text: " Project:  Project Name" -> " Project:  Cabin Boy Fowx"
# Only the final run ("Project Name") changes; re-apply formatting to the
# preceding runs so the engine doesn't collapse them into one plain run.
$c2 = $ws.Range("C2")
$c2Title = $c2.Characters(12, 12)
$c2Title.Text = "Cabin Boy Fowx"

$c2Run1 = $c2.Characters(1, 10)
$c2Run1.Font.Size = 26

$c2Run2 = $c2.Characters(11, 1)
$c2Run2.Font.Size = 36

$c2Run3 = $c2.Characters(12, 14)
$c2Run3.Font.Bold = $true
$c2Run3.Font.Size = 36
$c2Run3.Font.Color = 2315831

# D2 rich text: "Developer:  Firstname Lastname" -> "Developer:  Ritikh Prasad"
$d2 = $ws.Range("D2")
$d2Name = $d2.Characters(13, 18)
$d2Name.Text = "Ritikh Prasad"

$d2Name2 = $d2.Characters(13, 13)
$d2Name2.Font.Bold = $true
$d2Name2.Font.Size = 20
$d2Name2.Font.Color = 2315831

# ---------------------------------------------------------------------------
# Milestone rows
# ---------------------------------------------------------------------------

# Row 5
$ws.Range("C5").Value = "Build the first pirate ship called the 'Crystal Voyager'. "

# Row 6
$ws.Range("B6").Value = "Player controls"
$ws.Range("C6").Value = "Code the controls for the player character, i.e movement and cannon fire"
$ws.Range("F6").Value = 43922

# Row 7
$ws.Range("B7").Value = "Creating enemy AI"
$ws.Range("C7").Value = 'Code and attach enemy AI to a placeholder enemy to test and refine code before creating the enemy ship (''Ruffian")'

# Row 8
$ws.Range("B8").Value = "Create ocean and land"
$ws.Range("C8").Value = "Create placeholder ocean and land"

# Row 9
$ws.Range("B9").Value = "Refine code"
$ws.Range("C9").Value = "Create a more advanced but simple PlayerControl script"

# Row 10
$ws.Range("B10").Value = "Refine Crystal Voyager"
$ws.Range("C10").Value = "Finish creating the Crystal Voyager"

# Row 11
$ws.Range("B11").Value = "Create 'Kill Switch' code"
$ws.Range("C11").Value = "Create the code for the Kill Switch which is the ship's special attack"

# Row 12
$ws.Range("B12").Value = "Code Crystal Voyager's Kill Switch"
$ws.Range("C12").Value = "Create the code for 'Crystal Fury', the Crystal Voyager's special attack/Kill Switch"

# Row 13
$ws.Range("B13").Value = "Animate Crystal Fury"
$ws.Range("C13").Value = "Create the animation for the Crystal Fury Kill Switch"

# ---------------------------------------------------------------------------
# Selection / view state
# ---------------------------------------------------------------------------
$ws.Range("C13").Select()
